# Merge the "https" + "://" runs into a single "https://" run on the
# Surfshark download-link paragraph (slide 8, shape 2, paragraph 5),
# matching how PowerPoint re-serialises a run after the user edits
# text inside it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(5)

# The paragraph currently reads "https" + "://" + "surfshark.com/download"
# as three separate runs. Re-typing the "https://" portion as one
# contiguous range merges it into a single run, leaving
# "surfshark.com/download" as the second run.
$prefix = $para.Characters(1, 8)
$prefix.Text = "https://"
